$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions): update F4 901 -> 902, F6 42 -> 43
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 902
$wsExpo.Range("F6").Value = 43

# Sheet "全部类型" (all types): update F5 901 -> 902, F7 42 -> 43
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 902
$wsAll.Range("F7").Value = 43
